$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "004" -> "002" (keep as text, leading zero must be preserved)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").ClearFormats()

# M2 / N2: date-like text strings
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric metric updates
$ws.Range("O2").Value = 113090467.6
$ws.Range("P2").Value = 164.0143748403
$ws.Range("Q2").Value = 370581488.65
$ws.Range("R2").Value = 537.4519398336
$ws.Range("S2").Value = 56068181.04
$ws.Range("T2").Value = 81.315320883
$ws.Range("U2").Value = -30686089.54
$ws.Range("V2").Value = -44.503837494
$ws.Range("Y2").Value = 30710843.37
$ws.Range("Z2").Value = 44.5397378138
$ws.Range("AA2").Value = -16257798.78
$ws.Range("AB2").Value = -23.5785805804
$ws.Range("AC2").Value = 68951558.48999999

# AD2 previously held a number (CCE_ADD_RATIO); now cleared to empty, like its
# neighboring empty cells (W2, X2, AE2, ...).
$ws.Range("AD2").NumberFormat = "@"
$ws.Range("AD2").Value = ""
$ws.Range("AD2").ClearFormats()
